$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 3 (Cancer row), shifting rows 3-9 down to 4-10
$ws.Rows("3:3").Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown)

# Populate the new row with the PLACES population figures
$ws.Range("A3").Value = "PLACES population"
$ws.Range("B3").Value = 88929
$ws.Range("C3").Value = 0

# Match the formatting used by the other header cells in column A (bold, left/top aligned)
$cell = $ws.Range("A3")
$cell.Font.Bold = $true
$cell.HorizontalAlignment = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignLeft
$cell.VerticalAlignment = [Microsoft.Office.Interop.Excel.XlVAlign]::xlVAlignTop

# Distinguish this row with a left/right thin border (no top/bottom)
$cell.Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeLeft).LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlContinuous
$cell.Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeRight).LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlContinuous

# Update the selected cell to match the saved view state
[void]$ws.Range("A17").Select()
